# Updates cached market-price / profit figures on several Sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) pulled in by the scheduled
# market-data runner. Each sheet stores its data in a Table (Table_<JOB>)
# spanning A1:N141; columns H-N hold the computed price/profit figures
# that get refreshed here. A few rows also have their HQ-profit (N) or
# NQ-profit (M) cell cleared entirely where no HQ/NQ recipe price applies.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1579.6666
$ws.Cells.Item(2, 10).Value = 292.66666
$ws.Cells.Item(2, 12).Value = 292.66666
$ws.Cells.Item(2, 14).Value = -518.66666
$ws.Cells.Item(6, 8).Value = 8710.9
$ws.Cells.Item(6, 9).Value = 8710.9
$ws.Cells.Item(6, 11).Value = 26132.7
$ws.Cells.Item(6, 13).Value = -26020.7
$ws.Cells.Item(17, 8).Value = 845.2766
$ws.Cells.Item(17, 10).Value = 845.2766
$ws.Cells.Item(17, 12).Value = 2535.8298
$ws.Cells.Item(17, 14).Value = -2871.8298
$ws.Cells.Item(28, 8).Value = 1115.4286
$ws.Cells.Item(28, 9).Value = 841.8
$ws.Cells.Item(28, 11).Value = 841.8
$ws.Cells.Item(28, 13).Value = -356.8
$ws.Cells.Item(40, 8).Value = 4211.3184
$ws.Cells.Item(40, 9).Value = 3049.8
$ws.Cells.Item(40, 11).Value = 3049.8
$ws.Cells.Item(40, 13).Value = -2874.8
$ws.Cells.Item(86, 8).Value = 9855.421
$ws.Cells.Item(86, 9).Value = 10111.846
$ws.Cells.Item(86, 11).Value = 10111.846
$ws.Cells.Item(86, 13).Value = -8988.846
$ws.Cells.Item(89, 8).Value = 9855.421
$ws.Cells.Item(89, 9).Value = 10111.846
$ws.Cells.Item(89, 11).Value = 50559.23
$ws.Cells.Item(89, 13).Value = -44943.23
$ws.Cells.Item(94, 8).Value = 2369.6
$ws.Cells.Item(94, 9).Value = 2369.6
$ws.Cells.Item(94, 11).Value = 2369.6
$ws.Cells.Item(94, 13).Value = -1918.6
$ws.Cells.Item(97, 8).Value = 1184.6666
$ws.Cells.Item(97, 10).Value = 1184.6666
$ws.Cells.Item(97, 12).Value = 3553.9998
$ws.Cells.Item(97, 14).Value = -4545.9998
$ws.Cells.Item(132, 8).Value = 1666.8125
$ws.Cells.Item(132, 9).Value = 1619.2858
$ws.Cells.Item(132, 10).Value = 1999.5
$ws.Cells.Item(132, 11).Value = 4857.857400000001
$ws.Cells.Item(132, 12).Value = 5998.5
$ws.Cells.Item(132, 13).Value = -2327.857400000001
$ws.Cells.Item(132, 14).Value = -11058.5
$ws.Cells.Item(133, 8).Value = 78799.8
$ws.Cells.Item(133, 10).Value = 78799.8
$ws.Cells.Item(133, 12).Value = 78799.8
$ws.Cells.Item(133, 14).Value = -88919.8
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(137, 8).Value = 3540.074
$ws.Cells.Item(137, 9).Value = 1949.45
$ws.Cells.Item(137, 11).Value = 5848.35
$ws.Cells.Item(137, 13).Value = -3298.35
$ws.Cells.Item(138, 8).Value = 2188.0925
$ws.Cells.Item(138, 10).Value = 2644.2368
$ws.Cells.Item(138, 12).Value = 7932.7104
$ws.Cells.Item(138, 14).Value = -18212.7104

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9618337
$ws.Cells.Item(32, 9).Value = 10002251
$ws.Cells.Item(32, 11).Value = 10002251
$ws.Cells.Item(32, 13).Value = -10001964
$ws.Cells.Item(132, 8).Value = 5283.1333
$ws.Cells.Item(132, 9).Value = 2374.4783
$ws.Cells.Item(132, 10).Value = 14840.143
$ws.Cells.Item(132, 11).Value = 7123.4349
$ws.Cells.Item(132, 12).Value = 44520.429
$ws.Cells.Item(132, 13).Value = -4593.4349
$ws.Cells.Item(132, 14).Value = -49580.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4802.4546
$ws.Cells.Item(20, 9).Value = 4335.4443
$ws.Cells.Item(20, 11).Value = 4335.4443
$ws.Cells.Item(20, 13).Value = -4088.4443
$ws.Cells.Item(134, 8).Value = 44244.5
$ws.Cells.Item(134, 9).Value = 792.7
$ws.Cells.Item(134, 11).Value = 2378.1
$ws.Cells.Item(134, 13).Value = 156.8999999999996

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 482.4
$ws.Cells.Item(107, 9).Value = 379.375
$ws.Cells.Item(107, 10).Value = 894.5
$ws.Cells.Item(107, 11).Value = 379.375
$ws.Cells.Item(107, 12).Value = 894.5
$ws.Cells.Item(107, 13).Value = 1540.625
$ws.Cells.Item(107, 14).Value = -4734.5
$ws.Cells.Item(132, 8).Value = 1190.35
$ws.Cells.Item(132, 9).Value = 1200.3684
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 3601.1052
$ws.Cells.Item(132, 12).Value = 3000
$ws.Cells.Item(132, 13).Value = -1071.1052
$ws.Cells.Item(132, 14).Value = -8060
$ws.Cells.Item(134, 8).Value = 560512
$ws.Cells.Item(134, 9).Value = 1000839.2
$ws.Cells.Item(134, 10).Value = 10103
$ws.Cells.Item(134, 11).Value = 3002517.6
$ws.Cells.Item(134, 12).Value = 30309
$ws.Cells.Item(134, 13).Value = -2999982.6
$ws.Cells.Item(134, 14).Value = -35379
$ws.Cells.Item(135, 8).Value = 99166.5
$ws.Cells.Item(135, 10).Value = 99166.5
$ws.Cells.Item(135, 12).Value = 99166.5
$ws.Cells.Item(135, 14).Value = -109306.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(88, 8).Value = 3880
$ws.Cells.Item(88, 10).Value = 4000
$ws.Cells.Item(88, 12).Value = 12000
$ws.Cells.Item(88, 14).Value = -12856
$ws.Cells.Item(91, 8).Value = 3880
$ws.Cells.Item(91, 10).Value = 4000
$ws.Cells.Item(91, 12).Value = 12000
$ws.Cells.Item(91, 14).Value = -14964
$ws.Cells.Item(137, 8).Value = 6376.143
$ws.Cells.Item(137, 10).Value = 4926.6
$ws.Cells.Item(137, 12).Value = 14779.8
$ws.Cells.Item(137, 14).Value = -24979.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 8001
$ws.Cells.Item(12, 9).Value = 2001.5
$ws.Cells.Item(12, 11).Value = 2001.5
$ws.Cells.Item(12, 13).Value = -1861.5
$ws.Cells.Item(70, 8).Value = 5623.25
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 14).ClearContents()
$ws.Cells.Item(73, 8).Value = 5623.25
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 4950
$ws.Cells.Item(113, 10).Value = 4950
$ws.Cells.Item(113, 12).Value = 4950
$ws.Cells.Item(113, 14).Value = -9290

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(48, 8).Value = 15861.25
$ws.Cells.Item(48, 10).Value = 15861.25
$ws.Cells.Item(48, 12).Value = 15861.25
$ws.Cells.Item(48, 14).Value = -17183.25
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 13).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 113833
$ws.Cells.Item(2, 9).Value = 113833
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 113833
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -113721
$ws.Cells.Item(2, 14).ClearContents()
$ws.Cells.Item(46, 8).Value = 53331
$ws.Cells.Item(46, 10).Value = 53331
$ws.Cells.Item(46, 12).Value = 53331
$ws.Cells.Item(46, 14).Value = -53793
$ws.Cells.Item(132, 8).Value = 1651.7407
$ws.Cells.Item(132, 9).Value = 1481.5
$ws.Cells.Item(132, 10).Value = 2400.8
$ws.Cells.Item(132, 11).Value = 4444.5
$ws.Cells.Item(132, 12).Value = 7202.400000000001
$ws.Cells.Item(132, 13).Value = -1914.5
$ws.Cells.Item(132, 14).Value = -12262.4
$ws.Cells.Item(134, 8).Value = 53331
$ws.Cells.Item(134, 10).Value = 53331
$ws.Cells.Item(134, 12).Value = 159993
$ws.Cells.Item(134, 14).Value = -165063
